# Generate Report for Handback
#
# The localization status report is refreshed: the de-de language is now
# handed back and in sync with en-US, its handback timestamp is bumped,
# the zh-cn handback timestamp is bumped, and the stale "handback file is
# not the latest" error detail is cleared on both language sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: zh-cn / de-de status columns ---
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus

# --- zh-cn sheet ---
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("K2").Value = "2016-08-19 04:44:27"
$zhcn.Range("P2").Value = ""

# --- de-de sheet ---
$dede.Range("C2").Value = $newStatus
$dede.Range("K2").Value = "2016-08-19 04:44:34"
$dede.Range("P2").Value = ""

# --- Column widths: re-fit now that cell contents changed length ---
$overview.Columns.Item(5).AutoFit()
$overview.Columns.Item(6).AutoFit()
$zhcn.Columns.Item(3).AutoFit()
$zhcn.Columns.Item(16).AutoFit()
$dede.Columns.Item(3).AutoFit()
$dede.Columns.Item(16).AutoFit()
